# Update cryptocurrency price/volume snapshot values.
# Price cells (column D) are stored as plain text in the workbook (some
# contain multiple "." separators that aren't valid numbers, e.g.
# "60.949.94"). Excel's COM Value setter auto-converts anything that
# parses as a number (e.g. "571.04") into a numeric cell, which would
# change the cell's stored type. Prefixing with a leading apostrophe
# forces Excel to keep/enter the value as literal text, exactly like
# typing '0.999 into a cell, so every D-cell stays text-typed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.949.94"
$ws.Range("E2").Value = "  +0.13%  "

$ws.Range("D3").Value = "'3.384.06"
$ws.Range("E3").Value = "  -0.76%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'571.04"

$ws.Range("D6").Value = "'141.70"
$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -0.44%  "

$ws.Range("E9").Value = "  +1.21%  "

$ws.Range("E10").Value = "  -1.21%  "

$ws.Range("E11").Value = "  -0.32%  "

$ws.Range("D12").Value = "'3.963.30"
$ws.Range("E12").Value = "  -0.72%  "

$ws.Range("E13").Value = "  +1.95%  "

$ws.Range("D14").Value = "'27.68"
$ws.Range("E14").Value = "  -2.15%  "

$ws.Range("D15").Value = "'3.383.12"
$ws.Range("E15").Value = "  -0.64%  "

$ws.Range("E16").Value = "  -0.83%  "

$ws.Range("D17").Value = "'61.078.40"
$ws.Range("E17").Value = "  +0.26%  "

$ws.Range("D18").Value = "'6.09"
$ws.Range("E18").Value = "  -3.54%  "

$ws.Range("D19").Value = "'13.61"
$ws.Range("E19").Value = "  -5.10%  "

$ws.Range("D20").Value = "'8.93"
$ws.Range("E20").Value = "  -4.27%  "

$ws.Range("D21").Value = "'381.85"
$ws.Range("E21").Value = "  -2.80%  "

$ws.Range("D22").Value = "'74.86"
$ws.Range("E22").Value = "  +2.68%  "

$ws.Range("E23").Value = "  -2.66%  "

$ws.Range("E24").Value = "  -0.07%  "

$ws.Range("D26").Value = "'3.521.84"
$ws.Range("E26").Value = "  -0.86%  "

$ws.Range("E27").Value = "  +0.93%  "

$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  +0.02%  "

$ws.Range("E29").Value = "  -2.04%  "

$ws.Range("E30").Value = "  -0.38%  "

$ws.Range("E31").Value = "  -2.59%  "

$ws.Range("D32").Value = "'1.40"
$ws.Range("E32").Value = "  -2.78%  "

$ws.Range("E33").Value = "  -0.08%  "

$ws.Range("D34").Value = "'23.30"
$ws.Range("E34").Value = "  -2.45%  "

$ws.Range("D35").Value = "'6.96"
$ws.Range("E35").Value = "  -0.74%  "

$ws.Range("D36").Value = "'166.05"
$ws.Range("E36").Value = "  -0.95%  "

$ws.Range("D37").Value = "'3.414.36"
$ws.Range("E37").Value = "  -0.65%  "

$ws.Range("D38").Value = "'5.00"
$ws.Range("E38").Value = "  -2.64%  "

$ws.Range("E39").Value = "  -4.90%  "

$ws.Range("E40").Value = "  -2.08%  "

$ws.Range("D41").Value = "'27.00"
$ws.Range("E41").Value = "  +0.06%  "

$ws.Range("E42").Value = "  +0.15%  "

$ws.Range("E43").Value = "  -1.97%  "

$ws.Range("E44").Value = "  -2.79%  "

$ws.Range("E45").Value = "  -2.77%  "

$ws.Range("E46").Value = "  -0.50%  "

$ws.Range("D47").Value = "'2.455.11"
$ws.Range("E47").Value = "  -5.34%  "

$ws.Range("D48").Value = "'22.93"
$ws.Range("E48").Value = "  -0.01%  "

$ws.Range("D49").Value = "'6.72"
$ws.Range("E49").Value = "  -3.55%  "

$ws.Range("D50").Value = "'0.0265"
$ws.Range("E50").Value = "  +1.84%  "

$ws.Range("E51").Value = "  +6.96%  "
